$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new row of data (A9=15, B9=81), matching the style of the existing
# data rows (A2:B8) which use the centered style.
$ws.Range("A9").Value = 15
$ws.Range("B9").Value = 81
$ws.Range("A9:B9").HorizontalAlignment = -4108

# Update the current selection to F12 (was E6)
$ws.Range("F12").Select()
